# Apply the 2021 ("R") column to the 3.b.1 vaccine-coverage sheet.
# Mirrors the prior year's ("Q") column formatting for every new R cell,
# then fills in the 2021 figures (or leaves the cell blank where Q was blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rows: taller to fit the new column / wrapped header text.
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(2).RowHeight = 15

# Row 3 - thin separator row, empty cells with the thick-bottom-border style.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)

# Row 4 - year header.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# Row 5 - section header, empty value cell.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

# Row 6
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 88.796593100633856

# Row 7
$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 86.908583391486388

# Row 8
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 89.680106631122953

# Row 9
$ws.Range("Q9").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = 95.775910364145659

# Row 10
$ws.Range("Q10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = 96.517042279754136

# Row 11
$ws.Range("Q11").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Range("R11").Value = 90.311530128242666

# Row 12
$ws.Range("Q12").Copy()
$ws.Range("R12").PasteSpecial(-4122)
$ws.Range("R12").Value = 90.746324915190343

# Row 13
$ws.Range("Q13").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("R13").Value = 90.894107952204379

# Row 14
$ws.Range("Q14").Copy()
$ws.Range("R14").PasteSpecial(-4122)
$ws.Range("R14").Value = 81.065680730752504

# Row 15
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("R15").Value = 85.088888888888889

# Row 16 - subtitle row, empty value cell.
$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)

# Row 17
$ws.Range("Q17").Copy()
$ws.Range("R17").PasteSpecial(-4122)
$ws.Range("R17").Value = 93.37839883628321

# Row 18
$ws.Range("Q18").Copy()
$ws.Range("R18").PasteSpecial(-4122)
$ws.Range("R18").Value = 93.091416608513612

# Row 19
$ws.Range("Q19").Copy()
$ws.Range("R19").PasteSpecial(-4122)
$ws.Range("R19").Value = 94.815061646117954

# Row 20
$ws.Range("Q20").Copy()
$ws.Range("R20").PasteSpecial(-4122)
$ws.Range("R20").Value = 100.53781512605042

# Row 21
$ws.Range("Q21").Copy()
$ws.Range("R21").PasteSpecial(-4122)
$ws.Range("R21").Value = 100.33525796237662

# Row 22
$ws.Range("Q22").Copy()
$ws.Range("R22").PasteSpecial(-4122)
$ws.Range("R22").Value = 93.78989283832054

# Row 23
$ws.Range("Q23").Copy()
$ws.Range("R23").PasteSpecial(-4122)
$ws.Range("R23").Value = 95.401432340746325

# Row 24
$ws.Range("Q24").Copy()
$ws.Range("R24").PasteSpecial(-4122)
$ws.Range("R24").Value = 92.308748798242007

# Row 25
$ws.Range("Q25").Copy()
$ws.Range("R25").PasteSpecial(-4122)
$ws.Range("R25").Value = 89.338842975206617

# Row 26
$ws.Range("Q26").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("R26").Value = 87.955555555555549

# Row 27 - section header, empty value cell.
$ws.Range("Q27").Copy()
$ws.Range("R27").PasteSpecial(-4122)

# Row 28
$ws.Range("Q28").Copy()
$ws.Range("R28").PasteSpecial(-4122)
$ws.Range("R28").Value = 89.631204460036727

# Row 29
$ws.Range("Q29").Copy()
$ws.Range("R29").PasteSpecial(-4122)
$ws.Range("R29").Value = 89.204466154919743

# Row 30
$ws.Range("Q30").Copy()
$ws.Range("R30").PasteSpecial(-4122)
$ws.Range("R30").Value = 84.751749416861045

# Row 31
$ws.Range("Q31").Copy()
$ws.Range("R31").PasteSpecial(-4122)
$ws.Range("R31").Value = 96.201680672268907

# Row 32
$ws.Range("Q32").Copy()
$ws.Range("R32").PasteSpecial(-4122)
$ws.Range("R32").Value = 95.567144719687093

# Row 33
$ws.Range("Q33").Copy()
$ws.Range("R33").PasteSpecial(-4122)
$ws.Range("R33").Value = 91.330444457457389

# Row 34
$ws.Range("Q34").Copy()
$ws.Range("R34").PasteSpecial(-4122)
$ws.Range("R34").Value = 91.368262344515642

# Row 35
$ws.Range("Q35").Copy()
$ws.Range("R35").PasteSpecial(-4122)
$ws.Range("R35").Value = 92.345373803964662

# Row 36
$ws.Range("Q36").Copy()
$ws.Range("R36").PasteSpecial(-4122)
$ws.Range("R36").Value = 88.660287081339717

# Row 37
$ws.Range("Q37").Copy()
$ws.Range("R37").PasteSpecial(-4122)
$ws.Range("R37").Value = 84.944444444444443

# Row 38 - footer row, empty value cell.
$ws.Range("Q38").Copy()
$ws.Range("R38").PasteSpecial(-4122)

# Move the active-sheet selection from C2 to R3, matching the author's
# final cursor position after adding the new column.
$ws.Range("R3").Select()
